$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 417.45456
$ws.Range("J17").Value = 417.45456
$ws.Range("L17").Value = 1252.36368
$ws.Range("N17").Value = -1588.36368

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H40").Value = 335419.62
$ws.Range("I40").Value = 2151.9473
$ws.Range("J40").Value = 911063.8
$ws.Range("K40").Value = 2151.9473
$ws.Range("L40").Value = 911063.8
$ws.Range("M40").Value = -1976.9473
$ws.Range("N40").Value = -911413.8

$ws.Range("H70").Value = 2464.9
$ws.Range("J70").Value = 1749.4286
$ws.Range("L70").Value = 5248.2858
$ws.Range("N70").Value = -5788.2858

$ws.Range("H73").Value = 2464.9
$ws.Range("J73").Value = 1749.4286
$ws.Range("L73").Value = 5248.2858
$ws.Range("N73").Value = -7120.2858

$ws.Range("H99").Value = 376.23077
$ws.Range("I99").Value = 314
$ws.Range("J99").Value = 475.8
$ws.Range("K99").Value = 942
$ws.Range("L99").Value = 1427.4
$ws.Range("M99").Value = 556
$ws.Range("N99").Value = -4423.4

$ws.Range("H112").Value = 986.65576
$ws.Range("J112").Value = 994.76666
$ws.Range("L112").Value = 2984.29998
$ws.Range("N112").Value = -5200.29998

$ws.Range("H127").Value = 2584.2693
$ws.Range("I127").Value = 575
$ws.Range("J127").Value = 2949.5908
$ws.Range("K127").Value = 1725
$ws.Range("L127").Value = 8848.7724
$ws.Range("M127").Value = 3235
$ws.Range("N127").Value = -18768.7724

$ws.Range("H138").Value = 2177.3215
$ws.Range("I138").Value = 1159
$ws.Range("J138").Value = 3352.3076
$ws.Range("K138").Value = 3477
$ws.Range("L138").Value = 10056.9228
$ws.Range("M138").Value = 1663
$ws.Range("N138").Value = -20336.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18870930
$ws.Range("I32").Value = 2618.6938
$ws.Range("K32").Value = 2618.6938
$ws.Range("M32").Value = -2331.6938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 50
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 63
$ws.Range("N3").ClearContents()

$ws.Range("H99").Value = 38463372
$ws.Range("I99").Value = 62501430
$ws.Range("J99").Value = 2474
$ws.Range("K99").Value = 62501430
$ws.Range("L99").Value = 2474
$ws.Range("M99").Value = -62499932
$ws.Range("N99").Value = -5470

$ws.Range("H122").Value = 17858050
$ws.Range("I122").Value = 62500370
$ws.Range("J122").Value = 1122.8
$ws.Range("K122").Value = 187501110
$ws.Range("L122").Value = 3368.4
$ws.Range("M122").Value = -187498660
$ws.Range("N122").Value = -8268.4

$ws.Range("H126").Value = 38463372
$ws.Range("I126").Value = 62501430
$ws.Range("J126").Value = 2474
$ws.Range("K126").Value = 187504290
$ws.Range("L126").Value = 7422
$ws.Range("M126").Value = -187501820
$ws.Range("N126").Value = -12362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 60.166668
$ws.Range("J2").Value = 52.2
$ws.Range("L2").Value = 313.2
$ws.Range("N2").Value = -539.2

$ws.Range("H4").Value = 58959.65
$ws.Range("I4").Value = 58959.65
$ws.Range("K4").Value = 176878.95
$ws.Range("M4").Value = -176766.95

$ws.Range("H38").Value = 32.272728
$ws.Range("I38").Value = 25
$ws.Range("J38").Value = 51.666668
$ws.Range("K38").Value = 75
$ws.Range("L38").Value = 155.000004
$ws.Range("M38").Value = 272
$ws.Range("N38").Value = -849.000004

$ws.Range("H121").Value = 571.75
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 629
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 1887
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -4507

$ws.Range("H134").Value = 20002558
$ws.Range("I134").Value = 35715140
$ws.Range("J134").Value = 4726.273
$ws.Range("K134").Value = 107145420
$ws.Range("L134").Value = 14178.819
$ws.Range("M134").Value = -107140350
$ws.Range("N134").Value = -24318.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 21744312
$ws.Range("I122").Value = 50010420
$ws.Range("J122").Value = 1150.3846
$ws.Range("K122").Value = 150031260
$ws.Range("L122").Value = 3451.1538
$ws.Range("M122").Value = -150028810
$ws.Range("N122").Value = -8351.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1053.375
$ws.Range("I16").Value = 895.63635
$ws.Range("J16").Value = 1400.4
$ws.Range("K16").Value = 895.63635
$ws.Range("L16").Value = 1400.4
$ws.Range("M16").Value = -725.63635
$ws.Range("N16").Value = -1740.4

$ws.Range("H69").Value = 37265.332
$ws.Range("J69").Value = 37265.332
$ws.Range("L69").Value = 37265.332
$ws.Range("N69").Value = -38887.332

$ws.Range("H72").Value = 37265.332
$ws.Range("J72").Value = 37265.332
$ws.Range("L72").Value = 111795.996
$ws.Range("N72").Value = -119907.996

$ws.Range("H132").Value = 24848928
$ws.Range("I132").Value = 35715624
$ws.Range("J132").Value = 10771.286
$ws.Range("K132").Value = 107146872
$ws.Range("L132").Value = 32313.858
$ws.Range("M132").Value = -107144342
$ws.Range("N132").Value = -37373.858

$ws.Range("H136").Value = 64937050
$ws.Range("I136").Value = 42330304
$ws.Range("J136").Value = 166667400
$ws.Range("K136").Value = 126990912
$ws.Range("L136").Value = 500002200
$ws.Range("M136").Value = -126988362
$ws.Range("N136").Value = -500007300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17642.742
$ws.Range("I122").Value = 21842.709
$ws.Range("J122").Value = 3242.8572
$ws.Range("K122").Value = 65528.12699999999
$ws.Range("L122").Value = 9728.571599999999
$ws.Range("M122").Value = -63078.12699999999
$ws.Range("N122").Value = -14628.5716

$ws.Range("H126").Value = 680.1905
$ws.Range("I126").Value = 425.4737
$ws.Range("J126").Value = 3100
$ws.Range("K126").Value = 1276.4211
$ws.Range("L126").Value = 9300
$ws.Range("M126").Value = 1193.5789
$ws.Range("N126").Value = -14240
